$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "68.624.14"
Set-TextValue $ws "E2" "  +1.02%  "
Set-TextValue $ws "D3" "3.279.52"
Set-TextValue $ws "E3" "  -0.02%  "
Set-TextValue $ws "E4" "  +0.02%  "
Set-TextValue $ws "D5" "584.18"
Set-TextValue $ws "E5" "  +0.35%  "
Set-TextValue $ws "D6" "186.34"
Set-TextValue $ws "E6" "  +1.82%  "
Set-TextValue $ws "E7" "  +0.03%  "
Set-TextValue $ws "E8" "  -0.39%  "
Set-TextValue $ws "D9" "0.133"
Set-TextValue $ws "E9" "  -0.96%  "
Set-TextValue $ws "E10" "  -0.95%  "
Set-TextValue $ws "E11" "  +0.60%  "
Set-TextValue $ws "D12" "3.852.79"
Set-TextValue $ws "E12" "  +0.05%  "
Set-TextValue $ws "E13" "  -0.26%  "
Set-TextValue $ws "D14" "28.43"
Set-TextValue $ws "E14" "  -0.61%  "
Set-TextValue $ws "D15" "68.692.81"
Set-TextValue $ws "E15" "  +1.15%  "
Set-TextValue $ws "E16" "  +1.38%  "
Set-TextValue $ws "D17" "3.377.68"
Set-TextValue $ws "E17" "  +3.02%  "
Set-TextValue $ws "D18" "5.88"
Set-TextValue $ws "E18" "  +0.56%  "
Set-TextValue $ws "E19" "  +1.16%  "
Set-TextValue $ws "D20" "394.91"
Set-TextValue $ws "E20" "  +4.62%  "
Set-TextValue $ws "E21" "  +0.72%  "
Set-TextValue $ws "D22" "71.69"
Set-TextValue $ws "E22" "  +0.58%  "
Set-TextValue $ws "D24" "0.521"
Set-TextValue $ws "E24" "  +1.28%  "
Set-TextValue $ws "E25" "  +0.19%  "
Set-TextValue $ws "E26" "  +4.87%  "
Set-TextValue $ws "E27" "  +1.43%  "
Set-TextValue $ws "E28" "  +0.02%  "
Set-TextValue $ws "D29" "5.77"
Set-TextValue $ws "E29" "  +0.98%  "
Set-TextValue $ws "E30" "  +0.12%  "
Set-TextValue $ws "D31" "23.08"
Set-TextValue $ws "E31" "  +0.53%  "
Set-TextValue $ws "D32" "7.19"
Set-TextValue $ws "E32" "  +3.34%  "
Set-TextValue $ws "E33" "  +1.19%  "
Set-TextValue $ws "E34" "  +0.00%  "
Set-TextValue $ws "E35" "  -1.21%  "
Set-TextValue $ws "D36" "163.39"
Set-TextValue $ws "E36" "  +0.55%  "
Set-TextValue $ws "E37" "  +7.76%  "
Set-TextValue $ws "E38" "  -3.24%  "
Set-TextValue $ws "D39" "26.79"
Set-TextValue $ws "E39" "  -1.19%  "
Set-TextValue $ws "E40" "  -0.28%  "
Set-TextValue $ws "D41" "6.60"
Set-TextValue $ws "E41" "  -2.80%  "
Set-TextValue $ws "E42" "  -3.08%  "
Set-TextValue $ws "D43" "25.71"
Set-TextValue $ws "E43" "  -0.48%  "
Set-TextValue $ws "B44" "Hedera"
Set-TextValue $ws "C44" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D44" "0.0693"
Set-TextValue $ws "E44" "  +1.52%  "
Set-TextValue $ws "B45" "OKB"
Set-TextValue $ws "C45" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D45" "41.37"
Set-TextValue $ws "E45" "  +1.19%  "
Set-TextValue $ws "D46" "2.661.57"
Set-TextValue $ws "E46" "  -0.51%  "
Set-TextValue $ws "D47" "340.59"
Set-TextValue $ws "E47" "  -3.06%  "
Set-TextValue $ws "E48" "  -0.26%  "
Set-TextValue $ws "D49" "6.37"
Set-TextValue $ws "E49" "  +3.04%  "
Set-TextValue $ws "D50" "31.91"
Set-TextValue $ws "E50" "  +2.37%  "
Set-TextValue $ws "E51" "  -0.48%  "
